$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"-0.0003960461763199419"
$ws.Range("F2").Value = [double]"71.64391326904297"
$ws.Range("G2").Value = [double]"-4.470348358154297e-08"
$ws.Range("I2").Value = [double]"0.1439132690429688"
$ws.Range("J2").Value = [double]"-4.470348358154297e-08"
$ws.Range("K2").Value = [double]"0.0003960461763199419"
$ws.Range("L2").Value = [double]"0.1439138139971618"
$ws.Range("E3").Value = [double]"-0.0003960461763199419"
$ws.Range("F3").Value = [double]"-71.64391326904297"
$ws.Range("G3").Value = [double]"3.725290298461914e-08"
$ws.Range("H3").Value = [double]"-1.192092895507812e-07"
$ws.Range("I3").Value = [double]"-0.1439132690429688"
$ws.Range("J3").Value = [double]"3.725290298461914e-08"
$ws.Range("K3").Value = [double]"0.0003959269670303911"
$ws.Range("L3").Value = [double]"0.1439138136691489"
$ws.Range("E4").Value = [double]"118.8882751464844"
$ws.Range("G4").Value = [double]"394.0816650390625"
$ws.Range("J4").Value = [double]"8.2322998046875"
$ws.Range("K4").Value = [double]"-1.925765991210938"
$ws.Range("L4").Value = [double]"8.642320046760759"
$ws.Range("E5").Value = [double]"118.8882751464844"
$ws.Range("F5").Value = [double]"-328.3188781738281"
$ws.Range("G5").Value = [double]"391.394287109375"
$ws.Range("H5").Value = [double]"116.0313110351562"
$ws.Range("I5").Value = [double]"-6.818878173828125"
$ws.Range("J5").Value = [double]"5.544921875"
$ws.Range("K5").Value = [double]"-2.856964111328125"
$ws.Range("L5").Value = [double]"9.241509729626431"
$ws.Range("E6").Value = [double]"118.8882751464844"
$ws.Range("F6").Value = [double]"328.8450012207031"
$ws.Range("G6").Value = [double]"-388.3965148925781"
$ws.Range("H6").Value = [double]"114.7068786621094"
$ws.Range("I6").Value = [double]"7.345001220703125"
$ws.Range("J6").Value = [double]"-2.547149658203125"
$ws.Range("K6").Value = [double]"-4.181396484375"
$ws.Range("L6").Value = [double]"8.827292386284613"
$ws.Range("E7").Value = [double]"118.8882751464844"
$ws.Range("F7").Value = [double]"-317.7643432617188"
$ws.Range("G7").Value = [double]"-389.7033081054688"
$ws.Range("H7").Value = [double]"115.2067718505859"
$ws.Range("I7").Value = [double]"3.73565673828125"
$ws.Range("J7").Value = [double]"-3.85394287109375"
$ws.Range("K7").Value = [double]"-3.681503295898438"
$ws.Range("L7").Value = [double]"6.508569231223672"
$ws.Range("E8").Value = [double]"125.0787048339844"
$ws.Range("G8").Value = [double]"1.877199172973633"
$ws.Range("H8").Value = [double]"121.8229141235352"
$ws.Range("J8").Value = [double]"1.877199172973633"
$ws.Range("K8").Value = [double]"-3.255790710449219"
$ws.Range("L8").Value = [double]"4.368765972142449"
$ws.Range("E9").Value = [double]"125.0787048339844"
$ws.Range("F9").Value = [double]"-177.1546325683594"
$ws.Range("G9").Value = [double]"-1.132626533508301"
$ws.Range("I9").Value = [double]"-3.886062622070312"
$ws.Range("J9").Value = [double]"-1.132626533508301"
$ws.Range("K9").Value = [double]"-2.630622863769531"
$ws.Range("L9").Value = [double]"4.827473689047516"
$ws.Range("E10").Value = [double]"-6.126372814178467"
$ws.Range("F10").Value = [double]"-431.0640430772268"
$ws.Range("G10").Value = [double]"384.1620574407201"
$ws.Range("H10").Value = [double]"-7.180988569464773"
$ws.Range("I10").Value = [double]"-7.795473983374734"
$ws.Range("J10").Value = [double]"2.325448554001355"
$ws.Range("K10").Value = [double]"-1.054615755286306"
$ws.Range("L10").Value = [double]"8.203007984518694"
$ws.Range("E11").Value = [double]"-6.126372814178467"
$ws.Range("F11").Value = [double]"-503.8510173463313"
$ws.Range("G11").Value = [double]"384.1620574407201"
$ws.Range("H11").Value = [double]"-7.180988569464773"
$ws.Range("I11").Value = [double]"-9.08244825247931"
$ws.Range("J11").Value = [double]"2.325448554001355"
$ws.Range("K11").Value = [double]"-1.054615755286306"
$ws.Range("L11").Value = [double]"9.434553069836937"
$ws.Range("E12").Value = [double]"125.0787048339844"
$ws.Range("F12").Value = [double]"-503.8510173463313"
$ws.Range("G12").Value = [double]"-1.921457457338567"
$ws.Range("H12").Value = [double]"121.0626772711425"
$ws.Range("I12").Value = [double]"-9.082463085049824"
$ws.Range("J12").Value = [double]"-1.921449827944035"
$ws.Range("K12").Value = [double]"-4.016027562841884"
$ws.Range("L12").Value = [double]"10.11491880926904"
$ws.Range("E13").Value = [double]"-6.126372814178467"
$ws.Range("F13").Value = [double]"-426.0194288919388"
$ws.Range("G13").Value = [double]"-390.0377156731189"
$ws.Range("H13").Value = [double]"-6.178564048389262"
$ws.Range("I13").Value = [double]"-2.750859798086822"
$ws.Range("J13").Value = [double]"-8.201106786400146"
$ws.Range("K13").Value = [double]"-0.0521912342107953"
$ws.Range("L13").Value = [double]"8.650324044542913"
$ws.Range("E14").Value = [double]"-6.126372814178467"
$ws.Range("F14").Value = [double]"-503.8510173463313"
$ws.Range("G14").Value = [double]"-390.0377156731189"
$ws.Range("H14").Value = [double]"-6.178564048389262"
$ws.Range("I14").Value = [double]"-9.08244825247931"
$ws.Range("J14").Value = [double]"-8.201106786400146"
$ws.Range("K14").Value = [double]"-0.0521912342107953"
$ws.Range("L14").Value = [double]"12.23730945534317"
$ws.Range("E15").Value = [double]"-6.126372814178467"
$ws.Range("F15").Value = [double]"502.0972939341241"
$ws.Range("G15").Value = [double]"392.6730547882676"
$ws.Range("H15").Value = [double]"-10.44466239838772"
$ws.Range("I15").Value = [double]"7.328724840272173"
$ws.Range("J15").Value = [double]"10.83644590154887"
$ws.Range("K15").Value = [double]"-4.318289584209253"
$ws.Range("L15").Value = [double]"13.77629821449535"
$ws.Range("E16").Value = [double]"-6.126372814178467"
$ws.Range("F16").Value = [double]"419.441972831314"
$ws.Range("G16").Value = [double]"392.6730547882676"
$ws.Range("H16").Value = [double]"-10.44466239838772"
$ws.Range("I16").Value = [double]"-3.826596262537976"
$ws.Range("J16").Value = [double]"10.83644590154887"
$ws.Range("K16").Value = [double]"-4.318289584209253"
$ws.Range("L16").Value = [double]"12.27676763919376"
$ws.Range("E17").Value = [double]"125.0787048339844"
$ws.Range("F17").Value = [double]"502.0972939341241"
$ws.Range("G17").Value = [double]"9.550821580280385"
$ws.Range("H17").Value = [double]"118.5462349611945"
$ws.Range("I17").Value = [double]"7.328739672842687"
$ws.Range("J17").Value = [double]"9.550829209674916"
$ws.Range("K17").Value = [double]"-6.532469872789875"
$ws.Range("L17").Value = [double]"13.6967852587235"
$ws.Range("E18").Value = [double]"-6.126372814178467"
$ws.Range("F18").Value = [double]"502.0972939341241"
$ws.Range("G18").Value = [double]"-379.8566428835929"
$ws.Range("H18").Value = [double]"-6.303154396674159"
$ws.Range("I18").Value = [double]"7.328724840272173"
$ws.Range("J18").Value = [double]"1.979966003125867"
$ws.Range("K18").Value = [double]"-0.1767815824956926"
$ws.Range("L18").Value = [double]"7.593531779472994"
$ws.Range("E19").Value = [double]"-6.126372814178467"
$ws.Range("F19").Value = [double]"435.4744029775859"
$ws.Range("G19").Value = [double]"-379.8566428835929"
$ws.Range("H19").Value = [double]"-6.303154396674159"
$ws.Range("I19").Value = [double]"12.20583388373387"
$ws.Range("J19").Value = [double]"1.979966003125867"
$ws.Range("K19").Value = [double]"-0.1767815824956926"
$ws.Range("L19").Value = [double]"12.3666445691121"
